$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new Price (column D), new Volume(1h) (column E).
# $null means "unchanged - leave as is".
$updates = @(
    @{ Row = 2; D = '28.902.52'; E = '  -1.14%  ' },
    @{ Row = 3; D = '1.831.46'; E = '  -1.55%  ' },
    @{ Row = 4; D = '0.9991'; E = '  -0.10%  ' },
    @{ Row = 5; D = '244.87'; E = '  +1.00%  ' },
    @{ Row = 6; D = '0.6915'; E = '  -1.00%  ' },
    @{ Row = 7; D = '0.9998'; E = '  -0.05%  ' },
    @{ Row = 8; D = '0.07664'; E = '  -2.03%  ' },
    @{ Row = 9; D = '0.3037'; E = '  -2.63%  ' },
    @{ Row = 10; D = '23.28'; E = '  -3.05%  ' },
    @{ Row = 11; D = '0.07805'; E = '  +0.08%  ' },
    @{ Row = 12; D = '92.99'; E = '  +1.09%  ' },
    @{ Row = 13; D = '1.824.39'; E = '  -2.05%  ' },
    @{ Row = 14; D = '5.087'; E = '  -0.98%  ' },
    @{ Row = 15; D = $null; E = '  -1.79%  ' },
    @{ Row = 16; D = '6.531'; E = '  -0.48%  ' },
    @{ Row = 17; D = '0.000008232'; E = '  -3.23%  ' },
    @{ Row = 18; D = '28.924.72'; E = '  -1.16%  ' },
    @{ Row = 19; D = '240.06'; E = '  -3.22%  ' },
    @{ Row = 20; D = '2.074.84'; E = '  -1.89%  ' },
    @{ Row = 21; D = $null; E = '  -2.13%  ' },
    @{ Row = 22; D = '0.9998'; E = '  -0.07%  ' },
    @{ Row = 23; D = '7.454'; E = '  -1.61%  ' },
    @{ Row = 24; D = '0.9998'; E = '  -0.07%  ' },
    @{ Row = 25; D = '0.1496'; E = '  -2.93%  ' },
    @{ Row = 26; D = '158.25'; E = '  -1.46%  ' },
    @{ Row = 27; D = '8.730'; E = '  -2.08%  ' },
    @{ Row = 28; D = $null; E = '  -2.47%  ' },
    @{ Row = 29; D = '1.538'; E = '  -2.64%  ' },
    @{ Row = 30; D = '4.218'; E = '  -1.42%  ' },
    @{ Row = 31; D = '4.134'; E = '  -2.72%  ' },
    @{ Row = 32; D = '1.195'; E = '  -1.06%  ' },
    @{ Row = 33; D = $null; E = '  -2.40%  ' },
    @{ Row = 34; D = '0.7765'; E = '  +1.96%  ' },
    @{ Row = 35; D = '1.851'; E = '  -1.48%  ' },
    @{ Row = 36; D = '1.142'; E = '  -2.89%  ' },
    @{ Row = 37; D = '2.690'; E = '  -0.20%  ' },
    @{ Row = 38; D = '1.277.85'; E = '  +2.58%  ' },
    @{ Row = 39; D = $null; E = '  -0.10%  ' },
    @{ Row = 40; D = $null; E = '  -1.60%  ' },
    @{ Row = 41; D = '0.9550'; E = '  +5.84%  ' },
    @{ Row = 42; D = '6.132'; E = '  +3.90%  ' },
    @{ Row = 43; D = '106.79'; E = '  -3.03%  ' },
    @{ Row = 44; D = '0.9994'; E = '  -0.05%  ' },
    @{ Row = 45; D = '9.674'; E = '  +1.28%  ' },
    @{ Row = 46; D = $null; E = '  -1.29%  ' },
    @{ Row = 47; D = '0.5163'; E = '  -0.46%  ' },
    @{ Row = 48; D = '1.975.01'; E = $null },
    @{ Row = 49; D = '63.67'; E = '  -7.49%  ' },
    @{ Row = 50; D = '1.751'; E = '  -1.10%  ' },
    @{ Row = 51; D = '6.955'; E = '  -0.86%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Leading apostrophe forces text so numeric-looking strings
        # (e.g. "0.9991") are not coerced to a Double, matching the
        # original inline-string cell type.
        $ws.Range("D" + $u.Row).Value = "'" + $u.D
        $ws.Range("D" + $u.Row).Style = "Normal"
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = "'" + $u.E
        $ws.Range("E" + $u.Row).Style = "Normal"
    }
}
